# Apply updated crypto price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.162.92"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.074.51"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.81"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.678"
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.82"
$ws.Range("E7").Value = "  +9.52%  "
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "61.52"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  +8.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.108"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.32"
$ws.Range("E13").Value = "  +6.72%  "
$ws.Range("D14").Value = "2.378.37"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.54"
$ws.Range("E16").Value = "  +7.39%  "
$ws.Range("D17").Value = "2.077.97"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "37.187.69"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.75"
$ws.Range("E19").Value = "  +7.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "74.75"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").Value = "0.0₃0928"
$ws.Range("E21").Value = "  +10.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.49"
$ws.Range("E22").Value = "  +5.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.27"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("E26").Value = "  +14.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.89"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.35"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.127"
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("E31").Value = "  +6.75%  "
$ws.Range("E32").Value = "  +6.61%  "
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("E34").Value = "  +10.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0909"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("B38").Value = "Cronos"
$ws.Range("C38").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.117"
$ws.Range("E38").Value = "  +26.39%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.80"
$ws.Range("E42").Value = "  -3.16%  "
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.18"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.39"
$ws.Range("E45").Value = "  +6.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.83"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.60"
$ws.Range("E47").Value = "  +14.66%  "
$ws.Range("E48").Value = "  +8.68%  "
$ws.Range("D49").Value = "1.307.07"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.94"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("E51").Value = "  -1.17%  "
